$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 69, shifting existing rows 69-152 down to 70-153.
$ws.Rows.Item(69).Insert()

# Populate the newly inserted row 69 with the new data.
$ws.Cells.Item(69, 1).Value = 4
$ws.Cells.Item(69, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(69, 3).Value = "Los Lagos"
$ws.Cells.Item(69, 4).Value = 44904
$ws.Cells.Item(69, 5).Value = 10
$ws.Cells.Item(69, 6).Value = 100112052
$ws.Cells.Item(69, 7).Value = "Albahaca"
$ws.Cells.Item(69, 8).Value = "Sin especificar"
$ws.Cells.Item(69, 9).Value = "Primera"
$ws.Cells.Item(69, 10).Value = 90
$ws.Cells.Item(69, 11).Value = 8000
$ws.Cells.Item(69, 12).Value = 8000
$ws.Cells.Item(69, 13).Value = 8000
$ws.Cells.Item(69, 14).Value = "$/docena de matas"
$ws.Cells.Item(69, 15).Value = "Región Metropolitana"
$ws.Cells.Item(69, 16).Value = 1333
$ws.Cells.Item(69, 17).Value = 6
$ws.Cells.Item(69, 18).Value = "Hortaliza"
